$d = $word.ActiveDocument

# The hours table is the (only) table in the document, with columns
# "Aktivität" / "Stunden" and a trailing "Gesamt" (total) row.
$t = $d.Tables(1)

# --- 1. "Debugging" now took 6 hours instead of 7 ---------------------
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    if ($t.Cell($i, 1).Range.Text -like "*Debugging*") {
        $t.Cell($i, 2).Range.Text = "6"
        break
    }
}

# --- 2. New activity row "Dokumentation erstellen" (4h), inserted ----
#        right before the "Gesamt" row -------------------------------
$gesamtRowIndex = -1
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    if ($t.Cell($i, 1).Range.Text -like "*Gesamt*") {
        $gesamtRowIndex = $i
        break
    }
}

$newRow = $t.Rows.Add($t.Rows($gesamtRowIndex))
$t.Cell($gesamtRowIndex, 1).Range.Text = "Dokumentation erstellen"
$t.Cell($gesamtRowIndex, 2).Range.Text = "4"

# the "Gesamt" row moved one index further down
$gesamtRowIndex = $gesamtRowIndex + 1

# --- 3. Update the grand total (142 -> 145) ---------------------------
$t.Cell($gesamtRowIndex, 2).Range.Text = "145"

# --- 4. The embedded OLE preview object was re-linked on save; bump
#        its bookkeeping ObjectID to match. --------------------------
$d.Fields.Item(1).Update()
